$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.327869176864624
$ws.Range("B1").Value = 1.346596598625183
$ws.Range("C1").Value = 3.835701465606689
$ws.Range("D1").Value = 3.482260465621948
$ws.Range("E1").Value = 1.058396339416504
